$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7338071465492249
$ws.Range("B1").Value = 0.7081208229064941
$ws.Range("C1").Value = 0.5288563370704651
$ws.Range("D1").Value = 0.497530072927475
$ws.Range("E1").Value = 0.5220888257026672
